$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Un-hide rows 111 and 114 (Apple TV+ / Crunchyroll gift-temp rows).
# ---------------------------------------------------------------------------
$ws.Rows.Item(111).Hidden = $false
$ws.Rows.Item(114).Hidden = $false

# ---------------------------------------------------------------------------
# 2) Fix row 133 (Youtube Premium + Music (x2)): CODIGO should be DIG010, and
#    the description cell becomes bold like the other Netflix rows below it.
# ---------------------------------------------------------------------------
$ws.Range("A133").Value = "DIG010"
$ws.Range("B133").Font.Bold = $true

# ---------------------------------------------------------------------------
# 3) Row 134 becomes the "Netflix 4K" plan (renamed from the generic
#    "Netflix - 1 Dispositivo"); bold its description like the new rows.
# ---------------------------------------------------------------------------
$ws.Range("B134").Value = "Netflix 4K - 1 Dispositivo"
$ws.Range("B134").Font.Bold = $true

# ---------------------------------------------------------------------------
# 4) Append three new Netflix plan rows (135-137) to the table.
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:L137"))

$netflixImg = "https://media4.giphy.com/avatars/NetflixisaJoke/APDgNHsUtxbA.png"
$subText = "🤖 Sucripcion x 30 dias."

# Row 135: Netflix Básico
$ws.Range("A135").Value = "DIG010"
$ws.Range("B135").Value = "Netflix Básico"
$ws.Range("B135").Font.Bold = $true
$ws.Range("C135").Value = 10
$ws.Range("E135").Value = "STREAMING"
$ws.Range("F135").Formula = "=+Tabla1[[#This Row],[CODIGO]]"
$ws.Range("G135").Value = $subText
$ws.Range("H135").Formula = "=28.9*1.1"
$ws.Range("I135").Value = $netflixImg

# Row 136: Netflix Estándar
$ws.Range("A136").Value = "DIG010"
$ws.Range("B136").Value = "Netflix Estándar"
$ws.Range("B136").Font.Bold = $true
$ws.Range("C136").Value = 10
$ws.Range("E136").Value = "STREAMING"
$ws.Range("F136").Formula = "=+Tabla1[[#This Row],[CODIGO]]"
$ws.Range("G136").Value = $subText
$ws.Range("H136").Formula = "=40.9*1.1"
$ws.Range("I136").Value = $netflixImg

# Row 137: Netflix Premium
$ws.Range("A137").Value = "DIG010"
$ws.Range("B137").Value = "Netflix Premium"
$ws.Range("B137").Font.Bold = $true
$ws.Range("C137").Value = 10
$ws.Range("E137").Value = "STREAMING"
$ws.Range("F137").Formula = "=+Tabla1[[#This Row],[CODIGO]]"
$ws.Range("G137").Value = $subText
$ws.Range("H137").Formula = "=52.9*1.1"
$ws.Range("I137").Value = $netflixImg

# ---------------------------------------------------------------------------
# 5) Drop the STOCK (colId=2) filter criteria, keep only the STREAMING
#    (colId=4 / column E) filter, now spanning the expanded table range.
# ---------------------------------------------------------------------------
$lo.Range.AutoFilter(5, @("STREAMING"), [Microsoft.Office.Interop.Excel.XlAutoFilterOperator]::xlFilterValues)

# ---------------------------------------------------------------------------
# 6) Extend the conditional formatting on column C to the new last row.
# ---------------------------------------------------------------------------
$fcs = $ws.Range("C2:C134").FormatConditions
$fc = $fcs.Item(1)
$fc.ModifyAppliesToRange($ws.Range("C2:C137"))

# ---------------------------------------------------------------------------
# 7) Update the visible selection to match the new rows.
# ---------------------------------------------------------------------------
$ws.Range("I135:I137").Select()
